$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.114836812019348
$ws.Range("B1").Value = 2.25190258026123
$ws.Range("C1").Value = 10.33907032012939
$ws.Range("D1").Value = 1.642699241638184
$ws.Range("E1").Value = 1.292449116706848
